$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41, shifting existing rows 41-71 down to 42-72
$ws.Rows.Item(41).Insert()

# Fill the new row 41 with the new data record
$ws.Cells.Item(41, 1).Value = 10
$ws.Cells.Item(41, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(41, 3).Value = "La Araucanía"
$ws.Cells.Item(41, 4).Value = 44603
$ws.Cells.Item(41, 5).Value = 9
$ws.Cells.Item(41, 6).Value = 100112030
$ws.Cells.Item(41, 7).Value = "Poroto granado"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 65
$ws.Cells.Item(41, 11).Value = 25000
$ws.Cells.Item(41, 12).Value = 25000
$ws.Cells.Item(41, 13).Value = 25000
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(41, 16).Value = 1000
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
